# Update the author credit line on slide 1: merge the first three runs
# ("Asger B. Breinholm, Mathias " + "Brandgaard" + " and Rob ") into a
# single run with corrected text ("Asger Breinholm, Mathias Brændgaard and
# Rob "), leaving the trailing "Bertojo" run untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# The subtitle placeholder has two paragraphs: "Group 4" and the credits
# line. The credits line is paragraph 2.
$para = $tr.Paragraphs(2, 1)

# The first three runs of that paragraph together read
# "Asger B. Breinholm, Mathias Brandgaard and Rob " (47 characters).
# Replace just that span, preserving the final "Bertojo" run/formatting.
$creditsStart = $para.Characters(1, 47)
$creditsStart.Text = "Asger Breinholm, Mathias Brændgaard and Rob "
